$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.655.03'
$ws.Range('D3').Value = '3.163.04'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '528.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.39%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.536'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +14.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.30'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.439'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.68%  '
$ws.Range('E11').Value = '  +3.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.140'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Value = '3.710.05'
$ws.Range('E13').Value = '  +1.68%  '
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('E15').Value = '  +3.43%  '
$ws.Range('D16').Value = '58.697.26'
$ws.Range('E16').Value = '  +1.65%  '
$ws.Range('D17').Value = '3.162.66'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('E18').Value = '  +3.33%  '
$ws.Range('E19').Value = '  +2.26%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.77'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.10%  '
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('E23').Value = '  +4.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +14.10%  '
$ws.Range('D28').Value = '0.0₃0859'
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.50%  '
$ws.Range('E30').Value = '  +0.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.13'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.21%  '
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('E34').Value = '  +4.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '156.75'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.34'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.43%  '
$ws.Range('D37').Value = '2.694.33'
$ws.Range('E37').Value = '  +7.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.99'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0691'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.36%  '
$ws.Range('E41').Value = '  +6.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.721'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0290'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.62%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '3.205.92'
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('E47').Value = '  +13.04%  '
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.978'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.02'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.52%  '
$ws.Range('E51').Value = '  +0.74%  '
